# "remove column from alcohol data"
# The measurement sheet (Sheet1) had an extra column (M) of data that is no
# longer wanted; the following column (N) held the values that should now
# become the new, final column. Deleting column M's entire column shifts
# everything from N (and beyond, though nothing else was populated) one
# column to the left, so the old N values become the new M values and the
# sheet's used range shrinks from A1:N119 to A1:M119.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("M1").EntireColumn.Delete()

# Leave the selection on the (now last) column M, matching the post-edit
# cursor position recorded in the saved workbook.
$ws.Range("M1").Select() | Out-Null
